# Updates to bldgs/SYDEC - updated with Canada specific residential ratios
# from bldgs/BCEU (Start Year Distributed Electricity Capacity).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SYDEC")

# Update distributed-electricity-capacity figures (MW) on the SYDEC sheet
$ws.Range("D3").Value = 51.268617947930522
$ws.Range("D5").Value = 46.467633592094053

$ws.Range("B6").Value = 0.26888826085434137
$ws.Range("C6").Value = 0.39610850719915114
$ws.Range("D6").Value = 123.7324803310796

$ws.Range("B7").Value = 346.08981762389692
$ws.Range("C7").Value = 509.83676483403775
$ws.Range("D7").Value = 378.01342402715841

$ws.Range("D9").Value = 19.303470914060487

$ws.Range("D11").Value = 3.5578499616884471

# Move the active selection on the sheet (matches saved view state)
$ws.Activate()
$ws.Range("K12").Select()
